# angular-cli upgrade to 1.7.1
#
# The only content-level change that is reachable through the PowerPoint
# COM automation surface is the re-colouring of the "Favicon" slide title
# (slide 2, shape 1): its run gets an explicit blue solid fill
# (RGB 0,112,192 / srgbClr 0070C0).
#
# Deleting the run and re-inserting the text before colouring it also
# drops the stray trailing <a:endParaRPr> that PowerPoint leaves behind
# after a plain Text= assignment, matching the target markup more closely.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

[void]$tr.Delete()
[void]$tr.InsertAfter("Favicon")
$tr.Font.Color.RGB = 0xC07000

Write-Output "Favicon title recoloured"
